$d = $word.ActiveDocument

# 1. Remove the two image paragraphs (the ones whose Range contains an InlineShape).
#    Delete from last to first so indices of earlier paragraphs stay valid.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $p.Range.Delete()
    }
}

# 2. Text fixes (typos / OCR-style corruption introduced by this revision).
$d.Content.Find.Execute(
    "C At point Y, the ball starts to drop as no force is acting on it.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "| C At point Y, the ball starts to drop as no.force is acting on it.", 2)

$d.Content.Find.Execute(
    "(3) Gand Donly",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(3) Cand D only", 2)

$d.Content.Find.Execute(
    "(4) A,BandDonly",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(4)  A,BandDonly", 2)

$d.Content.Find.Execute(
    "(1) " + [char]0x2018 + "Landing only",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "{1}. _Lanaing only", 2)

$d.Content.Find.Execute(
    "(4) _ Taking off, flying andanding",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(4) _ Taking off, flying andtanding", 2)
